$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the analysis timestamp (shared string in A2)
$ws.Range("A2").Value = "2025-05-30 10:35:30"

# Update the metric values in row 2 that changed
$ws.Range("B2").Value = 16251
$ws.Range("C2").Value = 11783
$ws.Range("D2").Value = 72.5063073041659
$ws.Range("E2").Value = 2237
$ws.Range("F2").Value = 13.76530675035383
$ws.Range("H2").Value = 19.1126699895391
$ws.Range("I2").Value = 9569
$ws.Range("J2").Value = 58.88253030582733
$ws.Range("K2").Value = 3088353.07
$ws.Range("L2").Value = 3576
$ws.Range("M2").Value = 22.00479970463356
$ws.Range("N2").Value = 1184662.23
$ws.Range("P2").Value = 28.95821795581811
$ws.Range("Q2").Value = 502454.2
$ws.Range("R2").Value = 3586
$ws.Range("S2").Value = 22.06633437942281
$ws.Range("T2").Value = 3430
$ws.Range("U2").Value = 21.1063934527106
$ws.Range("V2").Value = 2430964.87
$ws.Range("W2").Value = 2298
$ws.Range("X2").Value = 14.14066826656821
$ws.Range("Y2").Value = 1433
$ws.Range("Z2").Value = 8.817918897298629
$ws.Range("AA2").Value = 154934
$ws.Range("AB2").Value = 798
$ws.Range("AC2").Value = 4.91046704818165
$ws.Range("AD2").Value = 16263
$ws.Range("AE2").Value = 10604
$ws.Range("AF2").Value = 65.20322203775441
$ws.Range("AG2").Value = 5659
$ws.Range("AH2").Value = 34.79677796224559
$ws.Range("AI2").Value = 486
$ws.Range("AJ2").Value = 864
$ws.Range("AK2").Value = 1338
$ws.Range("AL2").Value = 18.08035714285714
$ws.Range("AM2").Value = 32.14285714285715
$ws.Range("AN2").Value = 49.77678571428572
$ws.Range("AO2").Value = 1337552.96
$ws.Range("AP2").Value = 251248.74
$ws.Range("AQ2").Value = 83665.06999999999
$ws.Range("AR2").Value = 79.97486012831214
$ws.Range("AS2").Value = 15.02264466516127
$ws.Range("AT2").Value = 5.002495206526584
$ws.Range("AU2").Value = 50.77415167810124
$ws.Range("AV2").Value = 228.8204225352113
$ws.Range("AW2").Value = 558.6183719193427

Write-Host "done"